# Updates cryptos list price/volume figures (and swaps the
# EnergySwap / ImmutableX ranking rows) per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.725.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.87%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.462.38'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.10%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.05%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.58%  '

# Row 7: USDC
$ws.Range("E7").Value = '  +0.04%  '

# Row 9: Toncoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.62'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.65%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +1.95%  '

# Row 11: Cardano
$ws.Range("E11").Value = '  +3.62%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.054.66'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.16%  '

# Row 13: Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.80'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.88%  '

# Row 14: TRON
$ws.Range("E14").Value = '  +2.27%  '

# Row 15: WrappedEther
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.464.15'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.43%  '

# Row 16: ShibaInu
$ws.Range("E16").Value = '  +0.30%  '

# Row 17: WrappedBTC
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.799.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.92%  '

# Row 18: Polkadot
$ws.Range("E18").Value = '  +3.75%  '

# Row 19: Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.81%  '

# Row 20: Uniswap
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.63%  '

# Row 21: BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.83%  '

# Row 22: Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.564'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.52%  '

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '

# Row 24: Dai
$ws.Range("E24").Value = '  -0.14%  '

# Row 25: WrappedeETH
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.608.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.33%  '

# Row 26: PEPE
$ws.Range("E26").Value = '  +1.25%  '

# Row 27: Kaspa
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.179'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.91%  '

# Row 28: RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.61%  '

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = '  +0.15%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("E30").Value = '  +1.12%  '

# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.14'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.08%  '

# Row 32: USDe
$ws.Range("E32").Value = '  -0.03%  '

# Row 33: Fetch.AI
$ws.Range("E33").Value = '  -0.79%  '

# Row 34: EthereumClassic
$ws.Range("E34").Value = '  +1.93%  '

# Row 35: Aptos
$ws.Range("E35").Value = '  +2.72%  '

# Row 36: NEARProtocol
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.98%  '

# Row 37: ImmutableX
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.57'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.62%  '

# Row 38: EnergySwap
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +20.66%  '

# Row 39: Monero
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.59'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.07%  '

# Row 40: RenzoRestakedETH
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.500.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.23%  '

# Row 41: Hedera
$ws.Range("E41").Value = '  +0.51%  '

# Row 42: Mantle
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.798'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.21%  '

# Row 43: Filecoin
$ws.Range("E43").Value = '  +1.57%  '

# Row 44: OKB
$ws.Range("E44").Value = '  -0.53%  '

# Row 45: Stacks
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.14%  '

# Row 46: ONDO
$ws.Range("E46").Value = '  +2.34%  '

# Row 47: Maker
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.594.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.62%  '

# Row 48: InjectiveProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.37'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.35%  '

# Row 49: Cosmos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.63%  '

# Row 50: dogwifhat
$ws.Range("E50").Value = '  +8.68%  '

# Row 51: FirstDigitalUSD
$ws.Range("E51").Value = '  +0.00%  '
